# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Wed Feb 21 08:30:35 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'51.635.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -0.07%  "

$ws.Range("D3").Value2 = "'2.948.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  +1.37%  "

$ws.Range("D4").Value2 = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -0.23%  "

$ws.Range("D5").Value2 = "'359.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +1.28%  "

$ws.Range("D6").Value2 = "'105.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -3.36%  "

$ws.Range("D7").Value2 = "'0.548"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  -2.38%  "

$ws.Range("E8").Value2 = "  -0.02%  "

$ws.Range("D9").Value2 = "'0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  -4.26%  "

$ws.Range("D10").Value2 = "'37.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -3.77%  "

$ws.Range("E11").Value2 = "  +2.48%  "

$ws.Range("D12").Value2 = "'0.0847"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -2.49%  "

$ws.Range("D13").Value2 = "'18.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -3.23%  "

$ws.Range("D14").Value2 = "'3.405.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.87%  "

$ws.Range("D15").Value2 = "'7.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -3.79%  "

$ws.Range("D16").Value2 = "'2.940.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +1.05%  "

$ws.Range("D17").Value2 = "'0.979"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +0.22%  "

$ws.Range("D18").Value2 = "'51.537.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -0.29%  "

$ws.Range("D19").Value2 = "'3.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -0.62%  "

$ws.Range("D20").Value2 = "'7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -2.03%  "

$ws.Range("D21").Value2 = "'13.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -4.09%  "

$ws.Range("D22").Value2 = "0.0₃0956"
$ws.Range("E22").Value2 = "  -1.98%  "

$ws.Range("D23").Value2 = "'69.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -2.01%  "

$ws.Range("D24").Value2 = "'264.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -1.45%  "

$ws.Range("E25").Value2 = "  -3.61%  "

$ws.Range("E26").Value2 = "  -4.80%  "

$ws.Range("D27").Value2 = "'26.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -1.10%  "

$ws.Range("E28").Value2 = "  +0.11%  "

$ws.Range("B29").Value2 = "Filecoin"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value2 = "'7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -3.61%  "

$ws.Range("B30").Value2 = "Hedera"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value2 = "'0.109"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +3.09%  "

$ws.Range("D31").Value2 = "'6.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +4.67%  "

$ws.Range("D32").Value2 = "'10.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -3.41%  "

$ws.Range("E33").Value2 = "  +0.82%  "

$ws.Range("D34").Value2 = "'35.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -5.37%  "

$ws.Range("D35").Value2 = "'51.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -2.07%  "

$ws.Range("B36").Value2 = "VeChain"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value2 = "'0.0426"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -2.57%  "

$ws.Range("B37").Value2 = "FirstDigitalUSD"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value2 = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +0.11%  "

$ws.Range("D38").Value2 = "'2.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +5.46%  "

$ws.Range("D39").Value2 = "'3.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +0.30%  "

$ws.Range("D40").Value2 = "'17.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -4.99%  "

$ws.Range("D41").Value2 = "'1.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -4.36%  "

$ws.Range("B42").Value2 = "EnergySwap"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value2 = "'23.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +3.47%  "

$ws.Range("B43").Value2 = "Stellar"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value2 = "'0.115"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -3.39%  "

$ws.Range("D44").Value2 = "'120.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +1.65%  "

$ws.Range("D45").Value2 = "'2.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -1.47%  "

$ws.Range("D46").Value2 = "'2.089.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -1.40%  "

$ws.Range("D47").Value2 = "'3.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -5.58%  "

$ws.Range("D48").Value2 = "'2.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -7.11%  "

$ws.Range("D49").Value2 = "'3.231.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +0.83%  "

$ws.Range("E50").Value2 = "  -4.33%  "

$ws.Range("D51").Value2 = "'0.0317"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -4.05%  "
